# Insert a new row at position 393 (pushing existing rows 393-419 down to 394-420)
# and populate it with the new data record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(393).Insert()

$ws.Range("A393").Value = 8
$ws.Range("B393").Value = "Terminal La Palmera de La Serena"
$ws.Range("C393").Value = "Coquimbo"
$ws.Range("D393").Value = 44714
$ws.Range("E393").Value = 4
$ws.Range("F393").Value = 100114001
$ws.Range("G393").Value = "Papa"
$ws.Range("H393").Value = "Asterix"
$ws.Range("I393").Value = "1a nueva(o)"
$ws.Range("J393").Value = 2200
$ws.Range("K393").Value = 9000
$ws.Range("L393").Value = 10000
$ws.Range("M393").Value = 9500
$ws.Range("N393").Value = "$/saco 25 kilos"
$ws.Range("O393").Value = "Provincia de Melipilla"
$ws.Range("P393").Value = 380
$ws.Range("Q393").Value = 25
$ws.Range("R393").Value = "Hortaliza"
